$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B2:B190 to the value 1 (overwriting the previous averages)
$ws.Range("B2:B190").Value = 1
